$wb = $excel.ActiveWorkbook

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldTimestamp)) {
            $cell.Value = $val.Replace($oldTimestamp, $newTimestamp)
        }
    }
}
